$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.405.60"
$ws.Range("E2").Value = "  -2.19%  "
$ws.Range("D3").Value = "2.049.18"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'241.64"
$ws.Range("E5").Value = "  -2.86%  "
$ws.Range("D6").Value = "'0.664"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'53.70"
$ws.Range("E8").Value = "  -8.33%  "
$ws.Range("D9").Value = "'58.11"
$ws.Range("E9").Value = "  -3.77%  "
$ws.Range("E10").Value = "  -7.96%  "
$ws.Range("D11").Value = "'0.0745"
$ws.Range("E11").Value = "  -5.56%  "
$ws.Range("E12").Value = "  -2.83%  "
$ws.Range("D13").Value = "'0.894"
$ws.Range("E13").Value = "  -2.52%  "
$ws.Range("E14").Value = "  -7.97%  "
$ws.Range("D15").Value = "2.349.82"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").Value = "'5.32"
$ws.Range("E16").Value = "  -8.53%  "
$ws.Range("D17").Value = "2.040.42"
$ws.Range("E17").Value = "  -2.02%  "
$ws.Range("D18").Value = "36.386.96"
$ws.Range("E18").Value = "  -2.14%  "
$ws.Range("D19").Value = "'16.51"
$ws.Range("E19").Value = "  -12.61%  "
$ws.Range("D20").Value = "'71.70"
$ws.Range("E20").Value = "  -4.87%  "
$ws.Range("D21").Value = "0.0₃0850"
$ws.Range("E21").Value = "  -6.78%  "
$ws.Range("D22").Value = "'236.30"
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("D23").Value = "'5.22"
$ws.Range("E23").Value = "  -5.29%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  -5.47%  "
$ws.Range("D26").Value = "'9.22"
$ws.Range("E26").Value = "  -4.08%  "
$ws.Range("E27").Value = "  -4.82%  "
$ws.Range("D28").Value = "'162.34"
$ws.Range("E28").Value = "  -5.33%  "
$ws.Range("D29").Value = "'19.95"
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("E30").Value = "  -3.55%  "
$ws.Range("D31").Value = "'5.05"
$ws.Range("E31").Value = "  -8.67%  "
$ws.Range("D32").Value = "'1.16"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("D33").Value = "'4.48"
$ws.Range("E33").Value = "  -6.93%  "
$ws.Range("D34").Value = "'0.0588"
$ws.Range("E34").Value = "  -6.62%  "
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("D37").Value = "'0.0825"
$ws.Range("E37").Value = "  -6.52%  "
$ws.Range("D38").Value = "'2.16"
$ws.Range("E38").Value = "  -8.01%  "
$ws.Range("D39").Value = "'1.23"
$ws.Range("E39").Value = "  -8.57%  "
$ws.Range("D40").Value = "'4.80"
$ws.Range("E40").Value = "  -7.69%  "
$ws.Range("E41").Value = "  -6.20%  "
$ws.Range("E42").Value = "  -6.08%  "
$ws.Range("D43").Value = "'2.79"
$ws.Range("E43").Value = "  -10.78%  "
$ws.Range("D44").Value = "'92.78"
$ws.Range("E44").Value = "  -8.35%  "
$ws.Range("E45").Value = "  -12.55%  "
$ws.Range("D46").Value = "1.375.52"
$ws.Range("E46").Value = "  +5.14%  "
$ws.Range("D47").Value = "'15.58"
$ws.Range("E47").Value = "  -9.67%  "
$ws.Range("D48").Value = "'7.28"
$ws.Range("E48").Value = "  +5.30%  "
$ws.Range("D49").Value = "'2.83"
$ws.Range("E49").Value = "  -1.78%  "
$ws.Range("D50").Value = "'2.25"
$ws.Range("E50").Value = "  -7.37%  "
$ws.Range("D51").Value = "2.238.27"
$ws.Range("E51").Value = "  -0.44%  "
